# Generate Report for Handback
# ------------------------------------------------------------------
# This localization-status workbook tracks, per language (zh-cn / de-de),
# the handoff/handback lifecycle of each source file. This edit records
# that a.md and b.md have come back from localization ("Handback"):
#   * Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview sheet + both language
#     sheets).
#   * Each language sheet grows two new columns for the two tracked rows:
#       E = Latest Target File   (the localized file, a.md)
#       F = Latest Handback File (the handback .xlf package)
#     both rendered as hyperlinks, matching the existing A/C hyperlink
#     styling.
#   * Latest Handback DateTime (column G) is stamped with the real
#     handback time instead of the "0001-01-01 00:00:00" placeholder.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: both language columns for a.md / b.md flip status.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------
# Helper data per language sheet: table name, handback datetime stamp,
# and the hyperlink targets to restore/add (kept in the same shape as
# the ones already on the sheet).
# ---------------------------------------------------------------

# ===================== zh-cn =====================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Range("G2").Value = "2016-03-08 06:16:26"
$ws.Range("G3").Value = "2016-03-08 06:16:26"

# Rebuild the hyperlinks collection so the new E/F links land in the
# same row-major order Excel would naturally keep them in.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bb1919041a530152cd2553229f54e1cca522d386/e2e/a.md", "", "", "a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd1e6eb810480a52b5fa70733ba2d830b08a0f34/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$ws.Range("E2").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dd1e6eb810480a52b5fa70733ba2d830b08a0f34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.md", "", "", "a.md") | Out-Null
$ws.Range("F2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dd1e6eb810480a52b5fa70733ba2d830b08a0f34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bb1919041a530152cd2553229f54e1cca522d386/e2e/b.md", "", "", "b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd1e6eb810480a52b5fa70733ba2d830b08a0f34/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$ws.Range("E3").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dd1e6eb810480a52b5fa70733ba2d830b08a0f34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.md", "", "", "a.md") | Out-Null
$ws.Range("F3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dd1e6eb810480a52b5fa70733ba2d830b08a0f34/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bb1919041a530152cd2553229f54e1cca522d386/.localization-config", "", "", ".localization-config") | Out-Null

$ws.Range("E2").Style = $ws.Range("A2").Style
$ws.Range("F2").Style = $ws.Range("C2").Style
$ws.Range("E3").Style = $ws.Range("A3").Style
$ws.Range("F3").Style = $ws.Range("C3").Style

# ===================== de-de =====================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Range("G2").Value = "2016-03-08 06:16:31"
$ws.Range("G3").Value = "2016-03-08 06:16:31"

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bb1919041a530152cd2553229f54e1cca522d386/e2e/a.md", "", "", "a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5444e8a32ab17a1ec77370dd3b9f2319165b4a0d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$ws.Range("E2").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5444e8a32ab17a1ec77370dd3b9f2319165b4a0d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.md", "", "", "a.md") | Out-Null
$ws.Range("F2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5444e8a32ab17a1ec77370dd3b9f2319165b4a0d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bb1919041a530152cd2553229f54e1cca522d386/e2e/b.md", "", "", "b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5444e8a32ab17a1ec77370dd3b9f2319165b4a0d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$ws.Range("E3").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5444e8a32ab17a1ec77370dd3b9f2319165b4a0d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.md", "", "", "a.md") | Out-Null
$ws.Range("F3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5444e8a32ab17a1ec77370dd3b9f2319165b4a0d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bb1919041a530152cd2553229f54e1cca522d386/.localization-config", "", "", ".localization-config") | Out-Null

$ws.Range("E2").Style = $ws.Range("A2").Style
$ws.Range("F2").Style = $ws.Range("C2").Style
$ws.Range("E3").Style = $ws.Range("A3").Style
$ws.Range("F3").Style = $ws.Range("C3").Style
